$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.734577178955078
$ws.Range("B1").Value = 2.553204774856567
$ws.Range("C1").Value = 2.693010807037354
$ws.Range("D1").Value = 3.035745620727539
$ws.Range("E1").Value = 3.404922723770142
